# Commit: [ADDITIONAL SCRAPING] added code to scrape more data about a
# player's batting performance in a match, also updated the excel sheets
#
# 1. Insert a new "Player Info" worksheet as the first sheet, with the
#    player's ID / NAME / BATTING_HAND / BOWL_STYLE.
# 2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#    "ODI Batting" and "ODI Bowling" sheets, replacing the full scorecard
#    URL with just the bare match code.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Player Info" sheet before "ODI Batting" -------------
# NOTE: Add() re-uses/re-purposes the worksheet reference that is passed
# in as the "Before" sheet, so any handle obtained prior to this call
# becomes stale/renamed. Re-fetch the "ODI Batting"/"ODI Bowling" sheets
# by name *after* the new sheet has been inserted.
$beforeSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($beforeSheet)
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4746"
$playerInfo.Range("B2").Value = "Hussain Talat"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

# --- 2. Rename MATCH_CARD_LINK -> MATCH_CODE, store bare match code ------

# "ODI Batting" sheet: link lives in column D
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4238"

# "ODI Bowling" sheet: link lives in column B
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4238"
